# Correccion a Diebold Mariano y revision de Cap1
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Matriz_Resultados": some win/loss/tie entries changed from -1/1 to 0
# ---------------------------------------------------------------------------
$wsMatriz = $wb.Worksheets.Item("Matriz_Resultados")
$wsMatriz.Range("F2").Value = 0
$wsMatriz.Range("J2").Value = 0
$wsMatriz.Range("G4").Value = 0
$wsMatriz.Range("H4").Value = 0
$wsMatriz.Range("B6").Value = 0
$wsMatriz.Range("D7").Value = 0
$wsMatriz.Range("D8").Value = 0
$wsMatriz.Range("J9").Value = 0
$wsMatriz.Range("B10").Value = 0
$wsMatriz.Range("I10").Value = 0

# ---------------------------------------------------------------------------
# Sheet "P_valores": recomputed p-value matrix
# ---------------------------------------------------------------------------
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.0000006055425283424398
$wsP.Range("D2").Value = 0.0000009860043272968966
$wsP.Range("E2").Value = 0.00002489603957500641
$wsP.Range("F2").Value = 0.00290885226347859
$wsP.Range("G2").Value = 0.0000003357996287700615
$wsP.Range("H2").Value = 0.0000002217369408086256
$wsP.Range("I2").Value = 0.000002403415779728135
$wsP.Range("J2").Value = 0.003642738788639077

$wsP.Range("B3").Value = 0.0000006055425283424398
$wsP.Range("D3").Value = 0.000004232962749650682
$wsP.Range("E3").Value = 0.00000182779939272848
$wsP.Range("F3").Value = 0.0000005738411972266988
$wsP.Range("G3").Value = 0.0001812937221012945
$wsP.Range("H3").Value = 0.0002170466552542383
$wsP.Range("I3").Value = 0.01362008903117728
$wsP.Range("J3").Value = 0.0000007668756385470488

$wsP.Range("B4").Value = 0.0000009860043272968966
$wsP.Range("C4").Value = 0.000004232962749650682
$wsP.Range("E4").Value = 0.000002766109294949359
$wsP.Range("F4").Value = 0.0000009591209819692637
$wsP.Range("G4").Value = 0.003570321107812635
$wsP.Range("H4").Value = 0.002064290420534798
$wsP.Range("I4").Value = 0.0718188357005054
$wsP.Range("J4").Value = 0.000001351074071731517

$wsP.Range("B5").Value = 0.00002489603957500641
$wsP.Range("C5").Value = 0.00000182779939272848
$wsP.Range("D5").Value = 0.000002766109294949359
$wsP.Range("F5").Value = 0.00003828200494293554
$wsP.Range("G5").Value = 0.05506329268037713
$wsP.Range("H5").Value = 0.1080842082647553
$wsP.Range("I5").Value = 0.3456897100379164
$wsP.Range("J5").Value = 0.0002192116284642776

$wsP.Range("B6").Value = 0.00290885226347859
$wsP.Range("C6").Value = 0.0000005738411972266988
$wsP.Range("D6").Value = 0.0000009591209819692637
$wsP.Range("E6").Value = 0.00003828200494293554
$wsP.Range("G6").Value = 0.0000002054709458132464
$wsP.Range("H6").Value = 0.000000127633897895052
$wsP.Range("I6").Value = 0.000004484362988277368
$wsP.Range("J6").Value = 0.01562590435471978

$wsP.Range("B7").Value = 0.0000003357996287700615
$wsP.Range("C7").Value = 0.0001812937221012945
$wsP.Range("D7").Value = 0.003570321107812635
$wsP.Range("E7").Value = 0.05506329268037713
$wsP.Range("F7").Value = 0.0000002054709458132464
$wsP.Range("H7").Value = 0.3970610326266533
$wsP.Range("I7").Value = 0.8114390071716264
$wsP.Range("J7").Value = 0.0003191691109358796

$wsP.Range("B8").Value = 0.0000002217369408086256
$wsP.Range("C8").Value = 0.0002170466552542383
$wsP.Range("D8").Value = 0.002064290420534798
$wsP.Range("E8").Value = 0.1080842082647553
$wsP.Range("F8").Value = 0.000000127633897895052
$wsP.Range("G8").Value = 0.3970610326266533
$wsP.Range("I8").Value = 0.8699844764701392
$wsP.Range("J8").Value = 0.0006865340194259772

$wsP.Range("B9").Value = 0.000002403415779728135
$wsP.Range("C9").Value = 0.01362008903117728
$wsP.Range("D9").Value = 0.0718188357005054
$wsP.Range("E9").Value = 0.3456897100379164
$wsP.Range("F9").Value = 0.000004484362988277368
$wsP.Range("G9").Value = 0.8114390071716264
$wsP.Range("H9").Value = 0.8699844764701392
$wsP.Range("J9").Value = 0.004901163512162299

$wsP.Range("B10").Value = 0.003642738788639077
$wsP.Range("C10").Value = 0.0000007668756385470488
$wsP.Range("D10").Value = 0.000001351074071731517
$wsP.Range("E10").Value = 0.0002192116284642776
$wsP.Range("F10").Value = 0.01562590435471978
$wsP.Range("G10").Value = 0.0003191691109358796
$wsP.Range("H10").Value = 0.0006865340194259772
$wsP.Range("I10").Value = 0.004901163512162299

# ---------------------------------------------------------------------------
# Sheet "Estadisticos_DM": recomputed DM statistic matrix
# ---------------------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = 8.574276642007824
$wsE.Range("D2").Value = 8.227910348157646
$wsE.Range("E2").Value = 6.156868523869941
$wsE.Range("F2").Value = 3.598200079038175
$wsE.Range("G2").Value = 9.006721085153284
$wsE.Range("H2").Value = 9.320434151679192
$wsE.Range("I2").Value = 7.619558925530872
$wsE.Range("J2").Value = 3.484923958120196

$wsE.Range("B3").Value = -8.574276642007824
$wsE.Range("D3").Value = -7.248609694337011
$wsE.Range("E3").Value = -7.803220945632306
$wsE.Range("F3").Value = -8.613089787057703
$wsE.Range("G3").Value = -5.038008485487055
$wsE.Range("H3").Value = -4.941016771503238
$wsE.Range("I3").Value = -2.820470964692455
$wsE.Range("J3").Value = -8.405235048179513

$wsE.Range("B4").Value = -8.227910348157646
$wsE.Range("C4").Value = 7.248609694337011
$wsE.Range("E4").Value = -7.526359868937279
$wsE.Range("F4").Value = -8.247286836602939
$wsE.Range("G4").Value = -3.495026929169753
$wsE.Range("H4").Value = -3.771376159443729
$wsE.Range("I4").Value = -1.947481322648428
$wsE.Range("J4").Value = -8.009291042480017

$wsE.Range("B5").Value = -6.156868523869941
$wsE.Range("C5").Value = 7.803220945632306
$wsE.Range("D5").Value = 7.526359868937279
$wsE.Range("F5").Value = -5.905975162405412
$wsE.Range("G5").Value = 2.092780974508306
$wsE.Range("H5").Value = 1.716602901094829
$wsE.Range("I5").Value = 0.9758720499041124
$wsE.Range("J5").Value = -4.935687015447058

$wsE.Range("B6").Value = -3.598200079038175
$wsE.Range("C6").Value = 8.613089787057703
$wsE.Range("D6").Value = 8.247286836602939
$wsE.Range("E6").Value = 5.905975162405412
$wsE.Range("G6").Value = 9.378887277044534
$wsE.Range("H6").Value = 9.750441405157677
$wsE.Range("I6").Value = 7.211440637009393
$wsE.Range("J6").Value = 2.750621455075693

$wsE.Range("B7").Value = -9.006721085153284
$wsE.Range("C7").Value = 5.038008485487055
$wsE.Range("D7").Value = 3.495026929169753
$wsE.Range("E7").Value = -2.092780974508306
$wsE.Range("F7").Value = -9.378887277044534
$wsE.Range("H7").Value = -0.8736244497302728
$wsE.Range("I7").Value = -0.2431192562678579
$wsE.Range("J7").Value = -4.735316299259548

$wsE.Range("B8").Value = -9.320434151679192
$wsE.Range("C8").Value = 4.941016771503238
$wsE.Range("D8").Value = 3.771376159443729
$wsE.Range("E8").Value = -1.716602901094829
$wsE.Range("F8").Value = -9.750441405157677
$wsE.Range("G8").Value = 0.8736244497302728
$wsE.Range("I8").Value = 0.1667069019325835
$wsE.Range("J8").Value = -4.334264197936608

$wsE.Range("B9").Value = -7.619558925530872
$wsE.Range("C9").Value = 2.820470964692455
$wsE.Range("D9").Value = 1.947481322648428
$wsE.Range("E9").Value = -0.9758720499041124
$wsE.Range("F9").Value = -7.211440637009393
$wsE.Range("G9").Value = 0.2431192562678579
$wsE.Range("H9").Value = -0.1667069019325835
$wsE.Range("J9").Value = -3.335730120814543

$wsE.Range("B10").Value = -3.484923958120196
$wsE.Range("C10").Value = 8.405235048179513
$wsE.Range("D10").Value = 8.009291042480017
$wsE.Range("E10").Value = 4.935687015447058
$wsE.Range("F10").Value = -2.750621455075693
$wsE.Range("G10").Value = 4.735316299259548
$wsE.Range("H10").Value = 4.334264197936608
$wsE.Range("I10").Value = 3.335730120814543

# ---------------------------------------------------------------------------
# Sheet "Resumen": rows re-sorted / values updated
# ---------------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("Resumen")

# Row 3: LSPM
$wsR.Range("B3").Value = 4
$wsR.Range("D3").Value = 3
$wsR.Range("E3").Value = 50

# Row 4: was DeepAR -> now AV-MCPS
$wsR.Range("A4").Value = "AV-MCPS"
$wsR.Range("C4").Value = 1
$wsR.Range("D4").Value = 4
$wsR.Range("F4").Value = 2.270159067680547

# Row 6: MCPS
$wsR.Range("C6").Value = 1
$wsR.Range("D6").Value = 4

# Row 7: was AV-MCPS -> now DeepAR
$wsR.Range("A7").Value = "DeepAR"
$wsR.Range("B7").Value = 2
$wsR.Range("C7").Value = 0
$wsR.Range("D7").Value = 6
$wsR.Range("E7").Value = 25
$wsR.Range("F7").Value = 2.181297666943298

# Row 8: was AREPD -> now Block Bootstrapping
$wsR.Range("A8").Value = "Block Bootstrapping"
$wsR.Range("B8").Value = 0
$wsR.Range("D8").Value = 2
$wsR.Range("E8").Value = 0
$wsR.Range("F8").Value = 7.637559217939003

# Row 9: was EnCQR-LSTM -> now AREPD
$wsR.Range("A9").Value = "AREPD"
$wsR.Range("B9").Value = 0
$wsR.Range("D9").Value = 2
$wsR.Range("E9").Value = 0
$wsR.Range("F9").Value = 7.052496579093455

# Row 10: was Block Bootstrapping -> now EnCQR-LSTM
$wsR.Range("A10").Value = "EnCQR-LSTM"
$wsR.Range("C10").Value = 5
$wsR.Range("D10").Value = 3
$wsR.Range("F10").Value = 4.908057177504944
